$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on price cells whose new values would otherwise
# be auto-parsed as numbers by Excel (losing the original text-cell semantics).
# NumberFormat must be applied cell-by-cell: a comma-joined multi-area Range
# string only formats the first area.
foreach ($addr in @("D5", "D8", "D9", "D17", "D18", "D19", "D20", "D24", "D25", "D26", "D27", "D29", "D33", "D39", "D47", "D49", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "28.413.28"
$ws.Range("E2").Value = "  -0.28%  "

# Row 3
$ws.Range("D3").Value = "1.574.57"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "211.96"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "44.40"
$ws.Range("E8").Value = "  -4.03%  "

# Row 9
$ws.Range("D9").Value = "23.83"
$ws.Range("E9").Value = "  -0.93%  "

# Row 10
$ws.Range("E10").Value = "  -0.55%  "

# Row 11
$ws.Range("E11").Value = "  -0.45%  "

# Row 12
$ws.Range("E12").Value = "  +1.46%  "

# Row 13
$ws.Range("D13").Value = "1.798.15"
$ws.Range("E13").Value = "  +0.19%  "

# Row 14
$ws.Range("D14").Value = "1.569.66"
$ws.Range("E14").Value = "  -0.18%  "

# Row 15
$ws.Range("E15").Value = "  -0.04%  "

# Row 16
$ws.Range("D16").Value = "28.415.46"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
$ws.Range("D17").Value = "0.516"
$ws.Range("E17").Value = "  -0.88%  "

# Row 18
$ws.Range("D18").Value = "61.73"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19
$ws.Range("D19").Value = "228.53"
$ws.Range("E19").Value = "  +0.64%  "

# Row 20
$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  +0.91%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0686"
$ws.Range("E21").Value = "  -0.99%  "

# Row 22
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("E23").Value = "  +1.75%  "

# Row 24
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  -1.80%  "

# Row 25
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  -1.33%  "

# Row 26
$ws.Range("D26").Value = "151.04"
$ws.Range("E26").Value = "  +0.28%  "

# Row 27
$ws.Range("D27").Value = "14.94"
$ws.Range("E27").Value = "  -0.23%  "

# Row 29
$ws.Range("D29").Value = "6.37"
$ws.Range("E29").Value = "  -1.00%  "

# Row 31
$ws.Range("E31").Value = "  +3.52%  "

# Row 32
$ws.Range("E32").Value = "  -2.39%  "

# Row 33
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("E34").Value = "  -1.11%  "

# Row 35
$ws.Range("D35").Value = "1.382.96"
$ws.Range("E35").Value = "  -0.92%  "

# Row 36
$ws.Range("E36").Value = "  +4.35%  "

# Row 37
$ws.Range("E37").Value = "  -2.61%  "

# Row 38
$ws.Range("E38").Value = "  -0.41%  "

# Row 39
$ws.Range("D39").Value = "2.66"
$ws.Range("E39").Value = "  +1.83%  "

# Row 40
$ws.Range("E40").Value = "  -1.24%  "

# Row 41
$ws.Range("E41").Value = "  -2.41%  "

# Row 42
$ws.Range("E42").Value = "  +3.31%  "

# Row 43
$ws.Range("E43").Value = "  -0.04%  "

# Row 44
$ws.Range("E44").Value = "  -0.31%  "

# Row 45
$ws.Range("E45").Value = "  -0.84%  "

# Row 46
$ws.Range("E46").Value = "  -4.41%  "

# Row 47
$ws.Range("D47").Value = "62.44"
$ws.Range("E47").Value = "  -0.87%  "

# Row 48
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.711.27"
$ws.Range("E48").Value = "  +0.19%  "

# Row 49
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "0.919"
$ws.Range("E49").Value = "  -6.20%  "

# Row 50
$ws.Range("E50").Value = "  +1.77%  "

# Row 51
$ws.Range("D51").Value = "85.57"
$ws.Range("E51").Value = "  -0.55%  "
